$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two adjacent duplicate rows (row 40 and row 41)
# both holding "'Rhapsody on a Windy Night'" / 1950. Remove the second
# duplicate (row 41); everything below shifts up by one row.
$ws.Rows(41).Delete()

# Update the view to reflect where the analysis continued scrolling to.
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("A56").Select()
